# Fruta / hortaliza, semanal
# Refreshes the weekly pull: Fecha (D), Volumen (M), Precio minimo/maximo/
# promedio ponderado (N/O/P) and Precio $/Kg (S) are rotated across rows
# 2-20 (row 6 keeps its original values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  Year=2021; Month=7; Day=29; M=20; N=20000; O=20000; P=20000; S=1000}
    @{Row=3;  Year=2021; Month=8; Day=5;  M=45; N=20000; O=20000; P=20000; S=1000}
    @{Row=4;  Year=2021; Month=9; Day=9;  M=30; N=22000; O=22000; P=22000; S=1100}
    @{Row=5;  Year=2021; Month=4; Day=14; M=45; N=22000; O=22000; P=22000; S=1100}
    @{Row=7;  Year=2021; Month=7; Day=26; M=50; N=20000; O=20000; P=20000; S=1000}
    @{Row=8;  Year=2021; Month=4; Day=12; M=65; N=22000; O=22000; P=22000; S=1100}
    @{Row=9;  Year=2021; Month=6; Day=30; M=25; N=20000; O=20000; P=20000; S=1000}
    @{Row=10; Year=2021; Month=4; Day=19; M=20; N=22000; O=22000; P=22000; S=1100}
    @{Row=11; Year=2021; Month=4; Day=5;  M=70; N=25000; O=25000; P=25000; S=1250}
    @{Row=12; Year=2021; Month=4; Day=8;  M=25; N=25000; O=25000; P=25000; S=1250}
    @{Row=13; Year=2021; Month=7; Day=8;  M=36; N=20000; O=20000; P=20000; S=1000}
    @{Row=14; Year=2021; Month=6; Day=29; M=38; N=20000; O=20000; P=20000; S=1000}
    @{Row=15; Year=2021; Month=7; Day=12; M=20; N=20000; O=20000; P=20000; S=1000}
    @{Row=16; Year=2021; Month=4; Day=6;  M=30; N=25000; O=25000; P=25000; S=1250}
    @{Row=17; Year=2021; Month=7; Day=23; M=45; N=20000; O=20000; P=20000; S=1000}
    @{Row=18; Year=2021; Month=9; Day=6;  M=45; N=20000; O=20000; P=20000; S=1000}
    @{Row=19; Year=2021; Month=4; Day=15; M=38; N=22000; O=22000; P=22000; S=1100}
    @{Row=20; Year=2021; Month=4; Day=21; M=30; N=22000; O=22000; P=22000; S=1100}
)

foreach ($item in $data) {
    $fecha = Get-Date -Year $item.Year -Month $item.Month -Day $item.Day -Hour 0 -Minute 0 -Second 0
    $ws.Range("D" + $item.Row).Value = $fecha
    $ws.Range("M" + $item.Row).Value = $item.M
    $ws.Range("N" + $item.Row).Value = $item.N
    $ws.Range("O" + $item.Row).Value = $item.O
    $ws.Range("P" + $item.Row).Value = $item.P
    $ws.Range("S" + $item.Row).Value = $item.S
}
